$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row from column A (Beteckning), data starts at row 2.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 236 }

# Column C holds the "Förändrad" (last changed) date, stored as serial 45202 (2023-10-03).
# Update it to 45203 (2023-10-04) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
